$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in rows 7 and 8 with new API entries (order of first-use controls
# the shared-string table order, so write C7, C8, D8, D7, D4 in that order)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "POST"
$ws.Range("C7").Value = "/api2/create_TestRecords/"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "GET"
$ws.Range("C8").Value = "/api2/TestRecords/<int:pk>"
$ws.Range("D8").Value = "해당 코드의 검사 데이터 조회 (해당 코드는 환자와 연관 X)"

$ws.Range("D7").Value = "검사 데이터 추가 ( 템플릿 )"

# Update D4: "환자 추가" -> "환자 추가 ( 템플릿 )"
$ws.Range("D4").Value = "환자 추가 ( 템플릿 )"

# Widen column D (engine rounds ColumnWidth to whole pixels using MDW=7,
# so this lands on the nearest reachable value to the target 49.5625)
$ws.Columns.Item(4).ColumnWidth = 48.857142857142854

# Update sheet view: scroll so column B is leftmost, select D12
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D12").Select() | Out-Null
